$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer row appended below the existing data (row 22).
# Leading apostrophes force Chat id / Phone number (numeric-looking
# strings) to be stored as text, matching the rest of the sheet, and
# ClearFormats() strips the transient "quote prefix" cell style that
# Excel registers for that so the new row ends up unstyled like every
# other row already on the sheet.
$ws.Cells.Item(22, 1).Value = "'616525392"
$ws.Cells.Item(22, 2).Value = "Nurbek"
$ws.Cells.Item(22, 3).Value = "Boboyev"
$ws.Cells.Item(22, 4).Value = "'+998946696195"
$ws.Range("A22:D22").ClearFormats()
